# Remove the 6 "selegiline" rows from the Antidepressant Medications list.
# These correspond to rxcui 859186, 859190, 859193, 865206, 865210, 865214
# (ingredient = selegiline), located at rows 152-157 in the original sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A152:D157").EntireRow.Delete()
